$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5954290
$ws.Range("J17").Value = 6412154
$ws.Range("L17").Value = 19236462
$ws.Range("N17").Value = -19236798
$ws.Range("H41").Value = 626.5
$ws.Range("I41").Value = 198.27272
$ws.Range("J41").Value = 2196.6667
$ws.Range("K41").Value = 198.27272
$ws.Range("L41").Value = 2196.6667
$ws.Range("M41").Value = 241.72728
$ws.Range("N41").Value = -3076.6667
$ws.Range("H55").Value = 738.4286
$ws.Range("I55").Value = 836.6667
$ws.Range("J55").Value = 664.75
$ws.Range("K55").Value = 836.6667
$ws.Range("L55").Value = 664.75
$ws.Range("M55").Value = -622.6667
$ws.Range("N55").Value = -1092.75
$ws.Range("H61").Value = 149.5
$ws.Range("I61").Value = 149.5
$ws.Range("K61").Value = 448.5
$ws.Range("M61").Value = -276.5
$ws.Range("H74").Value = 16193.591
$ws.Range("I74").Value = 16459
$ws.Range("K74").Value = 16459
$ws.Range("M74").Value = -15523
$ws.Range("H77").Value = 16193.591
$ws.Range("I77").Value = 16459
$ws.Range("K77").Value = 82295
$ws.Range("M77").Value = -77615
$ws.Range("H82").Value = 4528
$ws.Range("I82").Value = 1038
$ws.Range("K82").Value = 3114
$ws.Range("M82").Value = -2708
$ws.Range("H85").Value = 4528
$ws.Range("I85").Value = 1038
$ws.Range("K85").Value = 3114
$ws.Range("M85").Value = -1710
$ws.Range("H96").Value = 695.125
$ws.Range("I96").Value = 287.16666
$ws.Range("J96").Value = 939.9
$ws.Range("K96").Value = 861.4999799999999
$ws.Range("L96").Value = 2819.7
$ws.Range("M96").Value = 511.5000200000001
$ws.Range("N96").Value = -5565.7
$ws.Range("H99").Value = 410.08334
$ws.Range("J99").Value = 1579
$ws.Range("L99").Value = 4737
$ws.Range("N99").Value = -7733
$ws.Range("H101").Value = 425
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H118").Value = 755.7143
$ws.Range("J118").Value = 616.1667
$ws.Range("L118").Value = 1848.5001
$ws.Range("N118").Value = -5162.5001
$ws.Range("H127").Value = 1080.5333
$ws.Range("I127").Value = 967.6667
$ws.Range("K127").Value = 2903.0001
$ws.Range("M127").Value = 2056.9999
$ws.Range("H131").Value = 1252072.9
$ws.Range("I131").Value = 2502296
$ws.Range("J131").Value = 1849.75
$ws.Range("K131").Value = 7506888
$ws.Range("L131").Value = 5549.25
$ws.Range("M131").Value = -7501848
$ws.Range("N131").Value = -15629.25
$ws.Range("H138").Value = 290476.25
$ws.Range("I138").Value = 3314.2
$ws.Range("J138").Value = 441614.2
$ws.Range("K138").Value = 9942.599999999999
$ws.Range("L138").Value = 1324842.6
$ws.Range("M138").Value = -4802.599999999999
$ws.Range("N138").Value = -1335122.6

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3462.8262
$ws.Range("I61").Value = 2320.625
$ws.Range("K61").Value = 2320.625
$ws.Range("M61").Value = -2108.625
$ws.Range("H74").Value = 219161.73
$ws.Range("J74").Value = 14907.125
$ws.Range("L74").Value = 14907.125
$ws.Range("N74").Value = -16655.125
$ws.Range("H77").Value = 219161.73
$ws.Range("J77").Value = 14907.125
$ws.Range("L77").Value = 74535.625
$ws.Range("N77").Value = -83271.625
$ws.Range("H123").Value = 90000
$ws.Range("J123").Value = 90000
$ws.Range("L123").Value = 90000
$ws.Range("N123").Value = -99800
$ws.Range("H136").Value = 3462.8262
$ws.Range("I136").Value = 2320.625
$ws.Range("K136").Value = 6961.875
$ws.Range("M136").Value = -4411.875

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 914.2222
$ws.Range("I22").Value = 914.2222
$ws.Range("K22").Value = 914.2222
$ws.Range("M22").Value = -741.2222
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H107").Value = 2850685.8
$ws.Range("I107").Value = 3847682.5
$ws.Range("K107").Value = 3847682.5
$ws.Range("M107").Value = -3845762.5
$ws.Range("H132").Value = 87799.8
$ws.Range("J132").Value = 87799.8
$ws.Range("L132").Value = 87799.8
$ws.Range("N132").Value = -97919.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1545.0212
$ws.Range("J16").Value = 1694
$ws.Range("L16").Value = 1694
$ws.Range("N16").Value = -2268
$ws.Range("H31").Value = 3493.0312
$ws.Range("I31").Value = 3013.08
$ws.Range("J31").Value = 5207.143
$ws.Range("K31").Value = 3013.08
$ws.Range("L31").Value = 5207.143
$ws.Range("M31").Value = -2718.08
$ws.Range("N31").Value = -5797.143
$ws.Range("H34").Value = 3493.0312
$ws.Range("I34").Value = 3013.08
$ws.Range("J34").Value = 5207.143
$ws.Range("K34").Value = 3013.08
$ws.Range("L34").Value = 5207.143
$ws.Range("M34").Value = -2811.08
$ws.Range("N34").Value = -5611.143
$ws.Range("H62").Value = 7701542.5
$ws.Range("I62").Value = 10004806
$ws.Range("J62").Value = 23998
$ws.Range("K62").Value = 10004806
$ws.Range("L62").Value = 23998
$ws.Range("M62").Value = -10004182
$ws.Range("N62").Value = -25246
$ws.Range("H65").Value = 7701542.5
$ws.Range("I65").Value = 10004806
$ws.Range("J65").Value = 23998
$ws.Range("K65").Value = 50024030
$ws.Range("L65").Value = 119990
$ws.Range("M65").Value = -50020910
$ws.Range("N65").Value = -126230
$ws.Range("H99").Value = 9499.267
$ws.Range("I99").Value = 10124.167
$ws.Range("J99").Value = 6999.6665
$ws.Range("K99").Value = 10124.167
$ws.Range("L99").Value = 6999.6665
$ws.Range("M99").Value = -8626.166999999999
$ws.Range("N99").Value = -9995.666499999999
$ws.Range("H113").Value = 1545.0212
$ws.Range("J113").Value = 1694
$ws.Range("L113").Value = 1694
$ws.Range("N113").Value = -6034
$ws.Range("H126").Value = 9499.267
$ws.Range("I126").Value = 10124.167
$ws.Range("J126").Value = 6999.6665
$ws.Range("K126").Value = 30372.501
$ws.Range("L126").Value = 20998.9995
$ws.Range("M126").Value = -27902.501
$ws.Range("N126").Value = -25938.9995

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 3857.7144
$ws.Range("I54").Value = 3004
$ws.Range("K54").Value = 9012
$ws.Range("M54").Value = -8453
$ws.Range("H103").Value = 1527.25
$ws.Range("J103").Value = 4990
$ws.Range("L103").Value = 14970
$ws.Range("N103").Value = -16728
$ws.Range("H113").Value = 6402.6665
$ws.Range("J113").Value = 7387.3335
$ws.Range("L113").Value = 22162.0005
$ws.Range("N113").Value = -26502.0005
$ws.Range("H132").Value = 2262.1714
$ws.Range("J132").Value = 3047.7896
$ws.Range("L132").Value = 27430.1064
$ws.Range("N132").Value = -32490.1064

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H36").Value = 4350
$ws.Range("J36").Value = 4350
$ws.Range("L36").Value = 4350
$ws.Range("N36").Value = -5320
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H43").Value = 47459.5
$ws.Range("J43").Value = 47459.5
$ws.Range("L43").Value = 47459.5
$ws.Range("N43").Value = -47761.5
$ws.Range("H113").Value = 3147.5945
$ws.Range("I113").Value = 3002.0605
$ws.Range("J113").Value = 4348.25
$ws.Range("K113").Value = 3002.0605
$ws.Range("L113").Value = 4348.25
$ws.Range("M113").Value = -832.0605
$ws.Range("N113").Value = -8688.25
$ws.Range("H122").Value = 8697.24
$ws.Range("I122").Value = 7955
$ws.Range("K122").Value = 23865
$ws.Range("M122").Value = -21415
$ws.Range("H132").Value = 1923.7894
$ws.Range("I132").Value = 1370.0667
$ws.Range("J132").Value = 4000.25
$ws.Range("K132").Value = 4110.2001
$ws.Range("L132").Value = 12000.75
$ws.Range("M132").Value = -1580.2001
$ws.Range("N132").Value = -17060.75
$ws.Range("H136").Value = 7522.6665
$ws.Range("J136").Value = 7522.6665
$ws.Range("L136").Value = 22567.9995
$ws.Range("N136").Value = -27667.9995

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H100").Value = 3692.4075
$ws.Range("I100").Value = 3633.6191
$ws.Range("J100").Value = 3898.1667
$ws.Range("K100").Value = 3633.6191
$ws.Range("L100").Value = 3898.1667
$ws.Range("M100").Value = -3092.6191
$ws.Range("N100").Value = -4980.1667
$ws.Range("H136").Value = 3128.8936
$ws.Range("I136").Value = 2860.7778
$ws.Range("K136").Value = 8582.3334
$ws.Range("M136").Value = -6032.3334

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 23811702
$ws.Range("I136").Value = 27027954
$ws.Range("K136").Value = 81083862
$ws.Range("M136").Value = -81081312
$ws.Range("H140").Value = 115255.29
$ws.Range("J140").Value = 115255.29
$ws.Range("L140").Value = 115255.29
$ws.Range("N140").Value = -125615.29
